$d = $word.ActiveDocument

# --- Change 1: "That was a new experience" -> "That was an experience" ---
# Restrict the find/replace range to the text after the en-dash run so the
# dash run itself is not touched/merged by the replace.
$r1 = $d.Range(5, 125)
$r1.Find.Execute("a new ", $true, $false, $false, $false, $false, $true, 1, $false, "an ", 2)

# Re-split the run right after the en-dash ("QB " | "-" | " That was an experience...")
# back into its own run (the Find/Replace above can coalesce adjacent runs
# with identical formatting). Adding then deleting a throwaway bookmark at
# that boundary forces Word to split the run there, and the split survives
# the bookmark's removal.
$tmp = $d.Range(4, 4)
$d.Bookmarks.Add("TempSplit", $tmp)
$d.Bookmarks.Item("TempSplit").Delete()

# Move the "_GoBack" bookmark (Word's "last edit" marker) to sit right after
# "That was an " and before "experience" -- this both relocates the bookmark
# and splits the run at that position, matching a real edit there.
$r2 = $d.Range(0, 0)
$r2.Find.Execute("experience, I am intrigued")
$editPos = $r2.Start
$bmRange = $d.Range($editPos, $editPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Change 2: drop the parenthetical stage direction after "1 - buzz " ---
# Match only the parenthetical's own run so the preceding "buzz " run is
# left completely untouched (avoids an unwanted run-merge).
$d.Content.Find.Execute("(maybe show drawn image of Nicholas Cage screaming BEES)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
